# Updates cryptocurrency Price (column D) and Volume(1h) (column E) cell
# values on Sheet1 to reflect the latest scrape, matching the commit
# "Updated cryptos list ... with GitHub Actions".
#
# All values in this sheet are stored as plain text (even when they look
# like numbers, e.g. "239.08" or "43.757.62"), so for any new value that
# Excel would otherwise auto-convert to a number we first force the cell's
# number format to Text ("@") and then assign the literal string. Values
# that Excel cannot parse as a number anyway (e.g. "43.748.18", which has
# two decimal points) are assigned directly.

$wb = $excel.ActiveWorkbook
if ($wb.Worksheets | Where-Object { $_.Name -eq "Sheet1" }) {
    $ws = $wb.Worksheets.Item("Sheet1")
} else {
    $ws = $wb.ActiveSheet
}

$ws.Range("D2").Value = '43.748.18'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '2.337.55'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.95'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.665'
$ws.Range("E6").Value = '  -4.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.93'
$ws.Range("E7").Value = '  -6.28%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.597'
$ws.Range("E9").Value = '  -5.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0992'
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.93'
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.20'
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("E14").Value = '  -5.57%  '
$ws.Range("D15").Value = '2.685.34'
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.04'
$ws.Range("E16").Value = '  -5.09%  '
$ws.Range("D18").Value = '2.337.35'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").Value = '43.654.70'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("E20").Value = '  -2.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '77.85'
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.57'
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.89'
$ws.Range("E23").Value = '  -1.97%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  +6.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.72'
$ws.Range("E26").Value = '  +2.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.48'
$ws.Range("E27").Value = '  -2.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.31'
$ws.Range("E28").Value = '  -6.67%  '
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '176.25'
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.13'
$ws.Range("E31").Value = '  -3.99%  '
$ws.Range("E32").Value = '  -2.13%  '
$ws.Range("E33").Value = '  -1.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0734'
$ws.Range("E34").Value = '  -3.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.05'
$ws.Range("E35").Value = '  -4.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.32'
$ws.Range("E36").Value = '  -1.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.71'
$ws.Range("E37").Value = '  -2.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.35'
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.35'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.74'
$ws.Range("E40").Value = '  +27.53%  '
$ws.Range("E41").Value = '  -3.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.84'
$ws.Range("E42").Value = '  +19.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.15'
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.107'
$ws.Range("E44").Value = '  +4.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.63'
$ws.Range("E45").Value = '  -3.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.195'
$ws.Range("E46").Value = '  -3.85%  '
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.22'
$ws.Range("E48").Value = '  -3.63%  '
$ws.Range("E49").Value = '  -5.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.60'
$ws.Range("E51").Value = '  -4.25%  '
